$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column C values per the data analysis refresh
$ws.Range("C2").Value = 0.01
$ws.Range("C7").Value = 0.4
$ws.Range("C11").Value = 0.3
$ws.Range("C14").Value = 0.2

# Update the active selection to match the new cursor position
$ws.Range("G15").Select()
